$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (bold, bordered, centered) from H1 so the
# new header cells I1/J1 reuse the same cell style as the other headers.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Set new header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new I/J numeric columns for rows 2-36
$values = @(
    @(1, 7),
    @(1, 6),
    @(1, 7),
    @(1, 6),
    @(1, 6),
    @(1, 7),
    @(1, 6),
    @(1, 6),
    @(1, 5),
    @(1, 7),
    @(1, 6),
    @(1, 5),
    @(11, 12),
    @(6, 7),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(6, 7),
    @(6, 8),
    @(7, 8),
    @(6, 8),
    @(8, 8),
    @(8, 8),
    @(3, 6),
    @(10, 10),
    @(8, 8),
    @(6, 7),
    @(7, 8),
    @(4, 4),
    @(5, 8),
    @(5, 8),
    @(1, 4),
    @(1, 3),
    @(4, 4),
    @(3, 4)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row = $row + 1
}
